# Apply updated Price (D) and Volume(1h) (E) figures to cryptos sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage so numeric-looking
# strings (e.g. "606.55") are not coerced into floating point numbers,
# then clear the temporary number-format override so no style is left behind.
function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range("D2").Value = '64.369.71'
$ws.Range("E2").Value = '  -2.93%  '

$ws.Range("D3").Value = '3.146.35'
$ws.Range("E3").Value = '  -1.83%  '

$ws.Range("E4").Value = '  -0.01%  '

Set-TextValue "D5" '606.55'
$ws.Range("E5").Value = '  -0.30%  '

Set-TextValue "D6" '149.87'
$ws.Range("E6").Value = '  -4.08%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").Value = '3.146.18'
$ws.Range("E8").Value = '  -1.89%  '

Set-TextValue "D9" '0.531'
$ws.Range("E9").Value = '  -3.48%  '

$ws.Range("E10").Value = '  -4.91%  '

Set-TextValue "D11" '5.59'
$ws.Range("E11").Value = '  -1.19%  '

Set-TextValue "D12" '0.479'
$ws.Range("E12").Value = '  -4.75%  '

$ws.Range("E13").Value = '  -3.31%  '

Set-TextValue "D14" '36.95'
$ws.Range("E14").Value = '  -3.88%  '

$ws.Range("D15").Value = '3.654.91'
$ws.Range("E15").Value = '  -2.11%  '

$ws.Range("D16").Value = '64.360.36'
$ws.Range("E16").Value = '  -3.15%  '

$ws.Range("E17").Value = '  +0.12%  '

$ws.Range("D18").Value = '3.138.61'
$ws.Range("E18").Value = '  -2.12%  '

$ws.Range("E19").Value = '  -4.04%  '

Set-TextValue "D20" '483.06'
$ws.Range("E20").Value = '  -4.61%  '

Set-TextValue "D21" '14.66'
$ws.Range("E21").Value = '  -4.27%  '

Set-TextValue "D22" '0.714'
$ws.Range("E22").Value = '  -2.25%  '

Set-TextValue "D23" '7.78'
$ws.Range("E23").Value = '  -2.59%  '

Set-TextValue "D24" '13.86'
$ws.Range("E24").Value = '  -5.28%  '

Set-TextValue "D25" '84.22'
$ws.Range("E25").Value = '  -1.13%  '

$ws.Range("E26").Value = '  +0.00%  '

Set-TextValue "D27" '2.94'
$ws.Range("E27").Value = '  -1.93%  '

Set-TextValue "D28" '8.60'
$ws.Range("E28").Value = '  -4.96%  '

$ws.Range("E31").Value = '  +1.44%  '

Set-TextValue "D32" '2.73'
$ws.Range("E32").Value = '  -6.70%  '

Set-TextValue "D33" '0.998'
$ws.Range("E33").Value = '  -0.34%  '

Set-TextValue "D34" '26.78'
$ws.Range("E34").Value = '  -5.22%  '

$ws.Range("E35").Value = '  -5.00%  '

$ws.Range("E36").Value = '  -5.10%  '

Set-TextValue "D37" '54.49'
$ws.Range("E37").Value = '  -1.66%  '

Set-TextValue "D38" '3.27'
$ws.Range("E38").Value = '  +7.46%  '

$ws.Range("D39").Value = '0.0₃0757'
$ws.Range("E39").Value = '  -1.32%  '

Set-TextValue "D40" '453.84'
$ws.Range("E40").Value = '  -9.58%  '

Set-TextValue "D41" '0.0402'
$ws.Range("E41").Value = '  -4.17%  '

$ws.Range("E42").Value = '  -5.20%  '

Set-TextValue "D43" '8.50'
$ws.Range("E43").Value = '  -2.42%  '

$ws.Range("D44").Value = '2.889.36'
$ws.Range("E44").Value = '  -0.82%  '

Set-TextValue "D45" '0.274'
$ws.Range("E45").Value = '  -7.58%  '

$ws.Range("E46").Value = '  -4.14%  '

Set-TextValue "D47" '26.79'
$ws.Range("E47").Value = '  -4.91%  '

$ws.Range("E48").Value = '  -0.06%  '

Set-TextValue "D49" '0.116'
$ws.Range("E49").Value = '  -1.00%  '

Set-TextValue "D50" '2.33'
$ws.Range("E50").Value = '  -3.24%  '

Set-TextValue "D51" '33.78'
$ws.Range("E51").Value = '  -0.07%  '

# Row 29 and 30 content swap (Hedera <-> ImmutableX) with new D/E values
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D29" '2.25'
$ws.Range("E29").Value = '  -4.25%  '

$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D30" '0.126'
$ws.Range("E30").Value = '  -1.82%  '
